# Weekly update: a new price-report row is inserted at the top of the
# data block (row 8 — right after the first six untouched entries),
# pushing every subsequent record down by one row. The oldest record
# that falls off the bottom of the original range reappears as the new
# last row (74).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8; this shifts rows 8:73 down to 9:74
# and grows the sheet's used range / dimension to A1:R74 automatically.
$ws.Rows("8:8").Insert()

# Populate the newly inserted row with this week's report.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 44649
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 100112022
$ws.Range("G8").Value = "Arveja Verde"
$ws.Range("H8").Value = "Perfection"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 30
$ws.Range("K8").Value = 25000
$ws.Range("L8").Value = 25000
$ws.Range("M8").Value = 25000
$ws.Range("N8").Value = "`$/saco 25 kilos"
$ws.Range("O8").Value = "Carahue"
$ws.Range("P8").Value = 1000
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
